$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old C column (full) and D column (tipo),
# shifting: old C(full)->E, old D(tipo)->F, old E(link)->G
$ws.Range("C:D").Insert()

# Update header row (row 1) for the new/shifted columns
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"
$ws.Range("E1").Value = "full"
$ws.Range("F1").Value = "tipo"
$ws.Range("G1").Value = "link"

$data = @(
    @{Row=2; Modelo='FONTE 90 BOB'; Politica='Acima'; Full='FULL'; Tipo='classico'; Link='https://www.mercadolivre.com.br/fonte-carregador-jfa-bob-storm-90a-bivolt-automatico-cor-preto/p/MLB21562641?pdp_filters=seller_id:84288564#searchVariation=MLB21562641&position=15&search_layout=grid&type=product&tracking_id=7908e0e6-88ac-44ae-9e00-6a9ea41713c5'},
    @{Row=3; Modelo='FONTE 120A'; Politica='Igual'; Full='FULL'; Tipo='classico'; Link='https://www.mercadolivre.com.br/fonte-automotiva-120a-amperes-jfa-carregador-cor-preto/p/MLB21392652?pdp_filters=seller_id:84288564#searchVariation=MLB21392652&position=32&search_layout=grid&type=product&tracking_id=7908e0e6-88ac-44ae-9e00-6a9ea41713c5'},
    @{Row=4; Modelo='FONTE 70A'; Politica='Igual'; Full='FULL'; Tipo='classico'; Link='https://www.mercadolivre.com.br/fonte-carregador-jfa-70a-bivolt-com-medidor-cca/p/MLB21455208?pdp_filters=seller_id:84288564#searchVariation=MLB21455208&position=39&search_layout=grid&type=product&tracking_id=7908e0e6-88ac-44ae-9e00-6a9ea41713c5'},
    @{Row=5; Modelo='FONTE 60A'; Politica='Acima'; Full='NA'; Tipo='premium'; Link='https://www.mercadolivre.com.br/fonte-carregador-jfa-60a-bivolt-storm-com-medidor-cca/p/MLB21320712?pdp_filters=seller_id:84288564#searchVariation=MLB21320712&position=26&search_layout=grid&type=product&tracking_id=6223e008-a5e3-4ee7-b710-65893a6434ad'},
    @{Row=6; Modelo='FONTE 200A'; Politica='Igual'; Full='NA'; Tipo='classico'; Link='https://www.mercadolivre.com.br/fonte-carregador-automotiva-jfa-200a-slim-bivolt-voltimetro/p/MLB21348561?pdp_filters=seller_id:84288564#searchVariation=MLB21348561&position=31&search_layout=grid&type=product&tracking_id=6223e008-a5e3-4ee7-b710-65893a6434ad'},
    @{Row=7; Modelo='FONTE 40A'; Politica='Acima'; Full='FULL'; Tipo='premium'; Link='https://www.mercadolivre.com.br/fonte-automotiva-40-amperes-jfa-storm-red-line-cca-sci-smart-cor-preto/p/MLB21621306?pdp_filters=seller_id:84288564#searchVariation=MLB21621306&position=17&search_layout=grid&type=product&tracking_id=0379eb80-5d24-4db3-a129-7475b014c4fd'},
    @{Row=8; Modelo='Sem Modelo'; Politica=$null; Full='NA'; Tipo='classico'; Link='https://www.mercadolivre.com.br/conversor-fio-para-rca-remoto-slim-12v-jfa-automotivo-cd-dvd/p/MLB25707531?pdp_filters=seller_id:84288564#searchVariation=MLB25707531&position=12&search_layout=grid&type=product&tracking_id=0379eb80-5d24-4db3-a129-7475b014c4fd'},
    @{Row=9; Modelo='Sem Modelo'; Politica=$null; Full='NA'; Tipo='premium'; Link='https://www.mercadolivre.com.br/controle-longa-distncia-reposico-jfa-tx-k1200-som/p/MLB29541981?pdp_filters=seller_id:84288564#searchVariation=MLB29541981&position=49&search_layout=grid&type=product&tracking_id=3814893e-ac89-40f2-9796-4fd8d0c1820a'},
    @{Row=10; Modelo='Sem Modelo'; Politica=$null; Full='FULL'; Tipo='classico'; Link='https://www.mercadolivre.com.br/controle-reposico-longa-distancia-avulso-jfa-k1200-preto/p/MLB28359785?pdp_filters=seller_id:84288564#searchVariation=MLB28359785&position=28&search_layout=grid&type=product&tracking_id=f0a0b00f-e747-4f48-8072-93aab8d8cdad'},
    @{Row=11; Modelo='Sem Modelo'; Politica=$null; Full='FULL'; Tipo='classico'; Link='https://produto.mercadolivre.com.br/MLB-777493114-filtro-anti-ruido-e-conversor-de-sinais-jfa-_JM#position%3D24%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D80f44a33-0d81-4b24-a1a7-200cbe8547c5'},
    @{Row=12; Modelo='Sem Modelo'; Politica=$null; Full='FULL'; Tipo='classico'; Link='https://produto.mercadolivre.com.br/MLB-2037736213-voltimetro-sequenciador-digital-jfa-vs5hi-led-vermelho-_JM#position%3D25%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D80f44a33-0d81-4b24-a1a7-200cbe8547c5'},
    @{Row=13; Modelo='Sem Modelo'; Politica=$null; Full='NA'; Tipo='classico'; Link='https://produto.mercadolivre.com.br/MLB-777493027-filtro-jfa-blindagem-e-anti-ruido-_JM#position%3D18%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3Dfcc1f06c-dce9-4daa-8c06-60176fe3e878'},
    @{Row=14; Modelo='Sem Modelo'; Politica=$null; Full='NA'; Tipo='premium'; Link='https://produto.mercadolivre.com.br/MLB-2037955677-conversor-interface-comando-remoto-rca-slim-jfa-_JM#position%3D22%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3Dfcc1f06c-dce9-4daa-8c06-60176fe3e878'},
    @{Row=15; Modelo='Sem Modelo'; Politica=$null; Full='NA'; Tipo='premium'; Link='https://produto.mercadolivre.com.br/MLB-3590118639-filtro-anti-ruido-jfa-eletromagnetico-2020k-rca-stereo-_JM#position%3D3%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D92832eb9-bf98-434a-91e0-9510d1df3908'},
    @{Row=16; Modelo='FONTE 120A'; Politica='Acima'; Full='FULL'; Tipo='premium'; Link='https://produto.mercadolivre.com.br/MLB-2641965235-fonte-jfa-storm-carregador-digital-automotivo-120a-bivolt-_JM#position%3D23%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D83808bd6-35f2-4c7c-b37b-20bcc28609ee'},
    @{Row=17; Modelo='FONTE 60A'; Politica='Acima'; Full='FULL'; Tipo='classico'; Link='https://produto.mercadolivre.com.br/MLB-2642396202-carregador-jfa-storm-fonte-60a-digital-medidor-de-cca-_JM#position%3D32%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D83808bd6-35f2-4c7c-b37b-20bcc28609ee'},
    @{Row=18; Modelo='Sem Modelo'; Politica=$null; Full='NA'; Tipo='classico'; Link='https://produto.mercadolivre.com.br/MLB-786094018-filtro-anti-ruido-jfa-rca-cd-dvd-eletromagnetico-_JM#position%3D28%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3Df818d2b9-60b3-42ee-b4ae-f2e174e1dfb1'},
    @{Row=19; Modelo='Sem Modelo'; Politica=$null; Full='NA'; Tipo='classico'; Link='https://produto.mercadolivre.com.br/MLB-3590131285-filtro-anti-ruido-jfa-2020k-rca-eletromagnetico-stereo-_JM#position%3D19%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D2ca30347-60c8-4bdc-9549-4f1b90e8352f'},
    @{Row=20; Modelo='Sem Modelo'; Politica=$null; Full='NA'; Tipo='premium'; Link='https://produto.mercadolivre.com.br/MLB-4571334324-conversor-interface-via-fio-para-rca-slim-jfa-comando-remoto-_JM#position%3D38%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D2ca30347-60c8-4bdc-9549-4f1b90e8352f'},
    @{Row=21; Modelo='Sem Modelo'; Politica=$null; Full='NA'; Tipo='premium'; Link='https://produto.mercadolivre.com.br/MLB-3453689301-controle-automotivo-longa-distncia-jfa-k1200-azul-scentral-_JM#position%3D22%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D10bd3ae5-60d3-4f6c-b4ca-2cca64fca9d5'},
    @{Row=22; Modelo='Sem Modelo'; Politica=$null; Full='NA'; Tipo='premium'; Link='https://produto.mercadolivre.com.br/MLB-777503537-filtro-blindagem-e-anti-ruido-jfa-eletromagnetico-_JM#position%3D5%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3Dcc843217-85ec-42c7-93b0-6149e9e65e47'},
    @{Row=23; Modelo='FONTE 60A'; Politica='Acima'; Full='NA'; Tipo='premium'; Link='https://produto.mercadolivre.com.br/MLB-2208861247-carregador-fonte-automotivo-60a-jfa-storm-com-display-_JM#position%3D24%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D826ed8e4-3daf-417c-8838-03c65c4ec462'},
    @{Row=24; Modelo='FONTE 120A'; Politica='Igual'; Full='NA'; Tipo='classico'; Link='https://produto.mercadolivre.com.br/MLB-2208906402-fonte-automotiva-carregador-120-amperes-bivolt-jfa-digital-_JM#position%3D25%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D826ed8e4-3daf-417c-8838-03c65c4ec462'},
    @{Row=25; Modelo='FONTE 120A'; Politica='Acima'; Full='NA'; Tipo='premium'; Link='https://produto.mercadolivre.com.br/MLB-2208906481-fonte-automotiva-120a-amperes-jfa-carregador-_JM#position%3D26%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D826ed8e4-3daf-417c-8838-03c65c4ec462'},
    @{Row=26; Modelo='Sem Modelo'; Politica=$null; Full='NA'; Tipo='classico'; Link='https://produto.mercadolivre.com.br/MLB-730728667-filtro-anto-ruido-jfa-eletromagnetica-_JM#position%3D17%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D9fe45666-f6d8-404e-a844-f35f919a571e'},
    @{Row=27; Modelo='FONTE 40A'; Politica='Igual'; Full='NA'; Tipo='premium'; Link='https://produto.mercadolivre.com.br/MLB-2208792205-fonte-carregador-40a-storm-jfa-medidor-cca-e-voltamp-_JM#position%3D23%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D9fe45666-f6d8-404e-a844-f35f919a571e'},
    @{Row=28; Modelo='FONTE 70A'; Politica='Igual'; Full='NA'; Tipo='premium'; Link='https://produto.mercadolivre.com.br/MLB-2208857188-carregador-fonte-som-70a-jfa-storm-digital-automotivo-bivolt-_JM#position%3D13%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D9c4e188a-8ee5-4928-a724-72e014ac0bb3'},
    @{Row=29; Modelo='Sem Modelo'; Politica=$null; Full='NA'; Tipo='premium'; Link='https://produto.mercadolivre.com.br/MLB-2037962261-voltimetro-sequenciador-jfa-vs5hi-digital-led-vermelho-_JM#position%3D36%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D9c4e188a-8ee5-4928-a724-72e014ac0bb3'},
    @{Row=30; Modelo='FONTE 40A'; Politica='Igual'; Full='NA'; Tipo='classico'; Link='https://produto.mercadolivre.com.br/MLB-2208785583-fonte-carregador-40a-storm-jfa-digital-com-medidor-de-cca-_JM#position%3D24%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D074be6ab-d92e-40c2-8d35-3c37104d8d7c'}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Modelo
    if ($null -eq $item.Politica) {
        # No "politica" value for this product (price band not applicable);
        # keep the cell present but empty, matching the default (unstyled) look.
        $ws.Cells.Item($r, 4).Style = $ws.Cells.Item($r, 1).Style
    } else {
        $ws.Cells.Item($r, 4).Value = $item.Politica
    }
    $ws.Cells.Item($r, 5).Value = $item.Full
    $ws.Cells.Item($r, 6).Value = $item.Tipo
    $ws.Cells.Item($r, 7).Value = $item.Link
}
